# naprawa problemow z rozgrywka
# Fix Uniki (Dodge) and Obrona (Defense) so they don't trigger too often / with too high values.

$wb = $excel.ActiveWorkbook

$wsWrogowie   = $wb.Worksheets.Item("Wrogowie")
$wsUzbrojenie = $wb.Worksheets.Item("Uzbrojenie")

# --- Wrogowie: reduce Bagiennik's starting agility (PoczatkowaZwinnosc), row 3, col E ---
$wsWrogowie.Range("E3").Value = 5

# --- Uzbrojenie: reduce the "Moc" (power) bonus of the defensive gear (Tarcza/Zbroja/Krucyfiks) ---
$wsUzbrojenie.Range("C2").Value = 2
$wsUzbrojenie.Range("C3").Value = 1
$wsUzbrojenie.Range("C4").Value = 2
$wsUzbrojenie.Range("C5").Value = 2

# --- Update selection / active-sheet UI state to match the author's final session ---
[void]$wsWrogowie.Activate()
[void]$wsWrogowie.Range("B2").Select()

[void]$wsUzbrojenie.Activate()
[void]$wsUzbrojenie.Range("C5").Select()
